$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.956883
$ws.Range("H2").Value = 35.913766
$ws.Range("I2").Value = 0.3392380274206944
$ws.Range("J2").Value = 0.2584869083704147
$ws.Range("M2").Value = 41.6173
$ws.Range("N2").Value = 83.2346
$ws.Range("O2").Value = 0.411761355892064
$ws.Range("P2").Value = 0.3324886731607734
$ws.Range("Q2").Value = 747.3169868758999
$ws.Range("R2").Value = 2989.2679475036
$ws.Range("S2").Value = 0.1396851101408943
$ws.Range("T2").Value = 0.08594396919350958
$ws.Range("G3").Value = 17.956883
$ws.Range("H3").Value = 35.913766
$ws.Range("I3").Value = 0.3392380274206944
$ws.Range("J3").Value = 0.2584869083704147
$ws.Range("O3").Value = 0.08885253351439082
$ws.Range("P3").Value = 0.1076198405427232
$ws.Range("Q3").Value = 161.2609019086053
$ws.Range("R3").Value = 967.5654114516319
$ws.Range("S3").Value = 0.03014215820075308
$ws.Range("T3").Value = 0.02781831986120553
$ws.Range("G4").Value = 17.956883
$ws.Range("H4").Value = 35.913766
$ws.Range("I4").Value = 0.3392380274206944
$ws.Range("J4").Value = 0.2584869083704147
$ws.Range("M4").Value = 13.45113833333333
$ws.Range("N4").Value = 40.353415
$ws.Range("O4").Value = 0.1330854946963174
$ws.Range("P4").Value = 0.1611956255073737
$ws.Range("Q4").Value = 241.5405172684816
$ws.Range("R4").Value = 1449.24310361089
$ws.Range("S4").Value = 0.04514766069908601
$ws.Range("T4").Value = 0.04166695888023619
$ws.Range("G5").Value = 17.956883
$ws.Range("H5").Value = 35.913766
$ws.Range("I5").Value = 0.3392380274206944
$ws.Range("J5").Value = 0.2584869083704147
$ws.Range("M5").Value = 11.258772
$ws.Range("N5").Value = 22.517544
$ws.Range("O5").Value = 0.1113942332731726
$ws.Range("P5").Value = 0.0899485109245354
$ws.Range("Q5").Value = 202.172451527676
$ws.Range("R5").Value = 808.6898061107039
$ws.Range("S5").Value = 0.03778915996163175
$ws.Range("T5").Value = 0.02325051250140562
$ws.Range("G6").Value = 17.956883
$ws.Range("H6").Value = 35.913766
$ws.Range("I6").Value = 0.3392380274206944
$ws.Range("J6").Value = 0.2584869083704147
$ws.Range("M6").Value = 15.450729
$ws.Range("N6").Value = 46.352187
$ws.Range("O6").Value = 0.1528694346476305
$ws.Range("P6").Value = 0.1851583014002596
$ws.Range("Q6").Value = 277.446932917707
$ws.Range("R6").Value = 1664.681597506242
$ws.Range("S6").Value = 0.05185912546277894
$ws.Range("T6").Value = 0.04786099688807053
$ws.Range("G7").Value = 17.956883
$ws.Range("H7").Value = 35.913766
$ws.Range("I7").Value = 0.3392380274206944
$ws.Range("J7").Value = 0.2584869083704147
$ws.Range("M7").Value = 10.313018
$ws.Range("N7").Value = 30.939054
$ws.Range("O7").Value = 0.1020369479764247
$ws.Range("P7").Value = 0.1235890484643348
$ws.Range("Q7").Value = 185.189657602894
$ws.Range("R7").Value = 1111.137945617364
$ws.Range("S7").Value = 0.03461481295555035
$ws.Range("T7").Value = 0.03194615104598723
$ws.Range("I8").Value = 0.6187742881378531
$ws.Range("J8").Value = 0.7072248972319991
$ws.Range("M8").Value = 41.6173
$ws.Range("N8").Value = 83.2346
$ws.Range("O8").Value = 0.411761355892064
$ws.Range("P8").Value = 0.3324886731607734
$ws.Range("Q8").Value = 1363.1152736129
$ws.Range("R8").Value = 8178.691641677399
$ws.Range("S8").Value = 0.2547873398747891
$ws.Range("T8").Value = 0.2351442677069317
$ws.Range("I9").Value = 0.6187742881378531
$ws.Range("J9").Value = 0.7072248972319991
$ws.Range("O9").Value = 0.08885253351439082
$ws.Range("P9").Value = 0.1076198405427232
$ws.Range("S9").Value = 0.05497966317461191
$ws.Range("T9").Value = 0.07611143066795155
$ws.Range("I10").Value = 0.6187742881378531
$ws.Range("J10").Value = 0.7072248972319991
$ws.Range("M10").Value = 13.45113833333333
$ws.Range("N10").Value = 40.353415
$ws.Range("O10").Value = 0.1330854946963174
$ws.Range("P10").Value = 0.1611956255073737
$ws.Range("Q10").Value = 440.5728413339316
$ws.Range("R10").Value = 3965.155572005385
$ws.Range("S10").Value = 0.08234988224218782
$ws.Range("T10").Value = 0.1140015596837002
$ws.Range("I11").Value = 0.6187742881378531
$ws.Range("J11").Value = 0.7072248972319991
$ws.Range("M11").Value = 11.258772
$ws.Range("N11").Value = 22.517544
$ws.Range("O11").Value = 0.1113942332731726
$ws.Range("P11").Value = 0.0899485109245354
$ws.Range("Q11").Value = 368.765010592356
$ws.Range("R11").Value = 2212.590063554136
$ws.Range("S11").Value = 0.0689278873962693
$ws.Range("T11").Value = 0.0636138263947759
$ws.Range("I12").Value = 0.6187742881378531
$ws.Range("J12").Value = 0.7072248972319991
$ws.Range("M12").Value = 15.450729
$ws.Range("N12").Value = 46.352187
$ws.Range("O12").Value = 0.1528694346476305
$ws.Range("P12").Value = 0.1851583014002596
$ws.Range("Q12").Value = 506.066580204717
$ws.Range("R12").Value = 4554.599221842453
$ws.Range("S12").Value = 0.09459167560212363
$ws.Range("T12").Value = 0.1309485606794501
$ws.Range("I13").Value = 0.6187742881378531
$ws.Range("J13").Value = 0.7072248972319991
$ws.Range("M13").Value = 10.313018
$ws.Range("N13").Value = 30.939054
$ws.Range("O13").Value = 0.1020369479764247
$ws.Range("P13").Value = 0.1235890484643348
$ws.Range("Q13").Value = 337.7881879133139
$ws.Range("R13").Value = 3040.093691219826
$ws.Range("S13").Value = 0.06313783984787137
$ws.Range("T13").Value = 0.0874052520991897
$ws.Range("G14").Value = 1.903653
$ws.Range("H14").Value = 3.807306
$ws.Range("I14").Value = 0.03596345137480081
$ws.Range("J14").Value = 0.02740282812891664
$ws.Range("M14").Value = 41.6173
$ws.Range("N14").Value = 83.2346
$ws.Range("O14").Value = 0.411761355892064
$ws.Range("P14").Value = 0.3324886731607734
$ws.Range("Q14").Value = 79.2248979969
$ws.Range("R14").Value = 316.8995919876
$ws.Range("S14").Value = 0.01480835950064629
$ws.Range("T14").Value = 0.009111129965436213
$ws.Range("G15").Value = 1.903653
$ws.Range("H15").Value = 3.807306
$ws.Range("I15").Value = 0.03596345137480081
$ws.Range("J15").Value = 0.02740282812891664
$ws.Range("O15").Value = 0.08885253351439082
$ws.Range("P15").Value = 0.1076198405427232
$ws.Range("Q15").Value = 17.095661852952
$ws.Range("R15").Value = 102.573971117712
$ws.Range("S15").Value = 0.003195443768572654
$ws.Range("T15").Value = 0.002949087993653659
$ws.Range("G16").Value = 1.903653
$ws.Range("H16").Value = 3.807306
$ws.Range("I16").Value = 0.03596345137480081
$ws.Range("J16").Value = 0.02740282812891664
$ws.Range("M16").Value = 13.45113833333333
$ws.Range("N16").Value = 40.353415
$ws.Range("O16").Value = 0.1330854946963174
$ws.Range("P16").Value = 0.1611956255073737
$ws.Range("Q16").Value = 25.606299841665
$ws.Range("R16").Value = 153.63779904999
$ws.Range("S16").Value = 0.004786213717202323
$ws.Range("T16").Value = 0.004417216020911774
$ws.Range("G17").Value = 1.903653
$ws.Range("H17").Value = 3.807306
$ws.Range("I17").Value = 0.03596345137480081
$ws.Range("J17").Value = 0.02740282812891664
$ws.Range("M17").Value = 11.258772
$ws.Range("N17").Value = 22.517544
$ws.Range("O17").Value = 0.1113942332731726
$ws.Range("P17").Value = 0.0899485109245354
$ws.Range("Q17").Value = 21.432795094116
$ws.Range("R17").Value = 85.73118037646401
$ws.Range("S17").Value = 0.00400612109175296
$ws.Range("T17").Value = 0.002464843585317024
$ws.Range("G18").Value = 1.903653
$ws.Range("H18").Value = 3.807306
$ws.Range("I18").Value = 0.03596345137480081
$ws.Range("J18").Value = 0.02740282812891664
$ws.Range("M18").Value = 15.450729
$ws.Range("N18").Value = 46.352187
$ws.Range("O18").Value = 0.1528694346476305
$ws.Range("P18").Value = 0.1851583014002596
$ws.Range("Q18").Value = 29.412826613037
$ws.Range("R18").Value = 176.476959678222
$ws.Range("S18").Value = 0.00549771247964335
$ws.Range("T18").Value = 0.00507386110991346
$ws.Range("G19").Value = 1.903653
$ws.Range("H19").Value = 3.807306
$ws.Range("I19").Value = 0.03596345137480081
$ws.Range("J19").Value = 0.02740282812891664
$ws.Range("M19").Value = 10.313018
$ws.Range("N19").Value = 30.939054
$ws.Range("O19").Value = 0.1020369479764247
$ws.Range("P19").Value = 0.1235890484643348
$ws.Range("Q19").Value = 19.632407654754
$ws.Range("R19").Value = 117.794445928524
$ws.Range("S19").Value = 0.003669600816983231
$ws.Range("T19").Value = 0.003386689453684514
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 0.6666666666666666
$ws.Range("G20").Value = 0.3188806666666666
$ws.Range("H20").Value = 0.956642
$ws.Range("I20").Value = 0.006024233066651711
$ws.Range("J20").Value = 0.006885366268669519
$ws.Range("M20").Value = 41.6173
$ws.Range("N20").Value = 83.2346
$ws.Range("O20").Value = 0.411761355892064
$ws.Range("P20").Value = 0.3324886731607734
$ws.Range("Q20").Value = 13.27095236886667
$ws.Range("R20").Value = 79.6257142132
$ws.Range("S20").Value = 0.002480546375734315
$ws.Range("T20").Value = 0.002289306294895874
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 0.6666666666666666
$ws.Range("G21").Value = 0.3188806666666666
$ws.Range("H21").Value = 0.956642
$ws.Range("I21").Value = 0.006024233066651711
$ws.Range("J21").Value = 0.006885366268669519
$ws.Range("O21").Value = 0.08885253351439082
$ws.Range("P21").Value = 0.1076198405427232
$ws.Range("Q21").Value = 2.863692095553778
$ws.Range("R21").Value = 25.773228859984
$ws.Range("S21").Value = 0.0005352683704531726
$ws.Range("T21").Value = 0.0007410020199124588
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 0.6666666666666666
$ws.Range("G22").Value = 0.3188806666666666
$ws.Range("H22").Value = 0.956642
$ws.Range("I22").Value = 0.006024233066651711
$ws.Range("J22").Value = 0.006885366268669519
$ws.Range("M22").Value = 13.45113833333333
$ws.Range("N22").Value = 40.353415
$ws.Range("O22").Value = 0.1330854946963174
$ws.Range("P22").Value = 0.1611956255073737
$ws.Range("Q22").Value = 4.289307959158888
$ws.Range("R22").Value = 38.60377163243
$ws.Range("S22").Value = 0.0008017380378412564
$ws.Range("T22").Value = 0.001109890922525555
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 0.6666666666666666
$ws.Range("G23").Value = 0.3188806666666666
$ws.Range("H23").Value = 0.956642
$ws.Range("I23").Value = 0.006024233066651711
$ws.Range("J23").Value = 0.006885366268669519
$ws.Range("M23").Value = 11.258772
$ws.Range("N23").Value = 22.517544
$ws.Range("O23").Value = 0.1113942332731726
$ws.Range("P23").Value = 0.0899485109245354
$ws.Range("Q23").Value = 3.590204721208
$ws.Range("R23").Value = 21.541228327248
$ws.Range("S23").Value = 0.0006710648235185605
$ws.Range("T23").Value = 0.0006193284430368478
$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 0.6666666666666666
$ws.Range("G24").Value = 0.3188806666666666
$ws.Range("H24").Value = 0.956642
$ws.Range("I24").Value = 0.006024233066651711
$ws.Range("J24").Value = 0.006885366268669519
$ws.Range("M24").Value = 15.450729
$ws.Range("N24").Value = 46.352187
$ws.Range("O24").Value = 0.1528694346476305
$ws.Range("P24").Value = 0.1851583014002596
$ws.Range("Q24").Value = 4.926938764006
$ws.Range("R24").Value = 44.342448876054
$ws.Range("S24").Value = 0.0009209211030846087
$ws.Range("T24").Value = 0.001274882722825492
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = 0.6666666666666666
$ws.Range("G25").Value = 0.3188806666666666
$ws.Range("H25").Value = 0.956642
$ws.Range("I25").Value = 0.006024233066651711
$ws.Range("J25").Value = 0.006885366268669519
$ws.Range("M25").Value = 10.313018
$ws.Range("N25").Value = 30.939054
$ws.Range("O25").Value = 0.1020369479764247
$ws.Range("P25").Value = 0.1235890484643348
$ws.Range("Q25").Value = 3.288622055185333
$ws.Range("R25").Value = 29.597598496668
$ws.Range("S25").Value = 0.0006146943560197984
$ws.Range("T25").Value = 0.000850955865473293
